$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "MuSCs" sending-cluster rows (old rows 6-9); rows 2-5 (FAPs sending
# cluster) shift up to become the only data rows, with the used range shrinking
# to A1:T5.
$ws.Range("A6:T9").EntireRow.Delete()

# Updated TPM-derived values for the remaining rows (FAPs -> ECs/FAPs/MuSCs/Resolving-Mac)
$ws.Range("I2:J2").Value = 1
$ws.Range("M2").Value = 1.716657
$ws.Range("N2").Value = 5.149971
$ws.Range("O2").Value = 0.3840886036988016
$ws.Range("P2").Value = 0.3840886036988015
$ws.Range("Q2").Value = 0.637614476196
$ws.Range("R2").Value = 5.738530285764
$ws.Range("S2").Value = 0.3840886036988016
$ws.Range("T2").Value = 0.3840886036988015

$ws.Range("I3:J3").Value = 1
$ws.Range("O3").Value = 0.07870146593648156
$ws.Range("P3").Value = 0.07870146593648154
$ws.Range("S3").Value = 0.07870146593648156
$ws.Range("T3").Value = 0.07870146593648154

$ws.Range("I4:J4").Value = 1
$ws.Range("M4").Value = 1.677572333333333
$ws.Range("N4").Value = 5.032717
$ws.Range("O4").Value = 0.3753437146230962
$ws.Range("P4").Value = 0.3753437146230962
$ws.Range("Q4").Value = 0.6230973366253334
$ws.Range("R4").Value = 5.607876029628
$ws.Range("S4").Value = 0.3753437146230962
$ws.Range("T4").Value = 0.3753437146230962

$ws.Range("I5:J5").Value = 1
$ws.Range("M5").Value = 0.7234496666666667
$ws.Range("N5").Value = 2.170349
$ws.Range("O5").Value = 0.1618662157416207
$ws.Range("P5").Value = 0.1618662157416207
$ws.Range("Q5").Value = 0.2687094627906667
$ws.Range("R5").Value = 2.418385165116
$ws.Range("S5").Value = 0.1618662157416207
$ws.Range("T5").Value = 0.1618662157416207
